$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "test3"
$ws.Range("B4").Value = "A team with that name already exists."

$ws.Range("A5").Value = "test4"
$ws.Range("B5").Value = "A team with that name already exists."

$ws.Range("A6").Value = "test5"
$ws.Range("B6").Value = "empty"

$ws.Range("D8").Select()
